# Updating filtered feeds from workflow
# Append one new feed-item row (row 96) to the "Filtered Feeds" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink  = "https://www.fiercebiotech.com/medtech/agilent-lands-fda-nod-companion-diagnostic-keytrudas-latest-cancer-approval"
$newTitle = '<a href="https://www.fiercebiotech.com/medtech/agilent-lands-fda-nod-companion-diagnostic-keytrudas-latest-cancer-approval" hreflang="en">Agilent companion diagnostic lands FDA nod alongside Keytruda''s latest cancer approval</a>'

$row = 96

$ws.Cells.Item($row, 1).Value = $newLink
$ws.Cells.Item($row, 2).Value = "companion diagnostic"
$ws.Cells.Item($row, 3).Value = $newTitle

$ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $newLink)
$ws.Cells.Item($row, 1).Style = "Hyperlink"
